# Rename worksheets: replace the space before "code" with a hyphen.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    $newName = $oldName -replace ' code', '-code'
    if ($newName -ne $oldName) {
        $ws.Name = $newName
    }
}
